# "Disconnected coke from pellet production" - remove the connections-sheet
# row describing coke being consumed (as fossil fuel) by the pellets /
# simple_pellets process, shifting the remaining rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Row 5 is: pellets | simple_pellets | inflow | fossil fuel -> coke | outflows | simple_coke | coke
$ws.Rows.Item(5).Delete()

# The edit leaves the "connections" tab active/selected (instead of "chains"),
# with the selection sitting on the row that is now last in the table.
$ws.Activate()
[void]$ws.Range("C19").Select()
